# Add a new "Hungary" sheet (FC Gallery Sounder FIM Panel test data) to the
# workbook, placed after "Spain" (the current last tab), cloned from the
# "Spain" sheet so that all formatting / styles / merged cells are preserved,
# then update the market name and part number for Hungary.

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# Spain was the previously-selected tab with a single-cell selection; leaving
# it mimics the "deselect" full-sheet selection state it is left in once the
# new sheet becomes active.
$spain.Cells.Select()

# Duplicate the "Spain" worksheet right after itself - this copies over the
# column widths, cell styles, merged cells and dimension automatically.
$spain.Copy([System.Reflection.Missing]::Value, $spain)

# The copy is inserted immediately after "Spain".
$hungary = $wb.Worksheets.Item($spain.Index + 1)
$hungary.Name = "Hungary"

# Update the market name and the part number for the Hungary test data.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3630/T3638/T3649"

# Make the new sheet the active tab with the same selection state as the
# original source data (active cell A9).
$hungary.Activate()
$hungary.Range("A9").Select()
